$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 15:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1432086
$ws.Range("C4").Value = 1738
$ws.Range("D4").Value = 310383
$ws.Range("E4").Value = 1036435
$ws.Range("F4").Value = 16337
$ws.Range("G4").Value = 71
$ws.Range("H4").Value = 85268

# Row 7 - Reino Unido
$ws.Range("B7").Value = 233151
$ws.Range("C7").Value = 3446
$ws.Range("E7").Value = 199193
$ws.Range("G7").Value = 428
$ws.Range("H7").Value = 33614

# Row 80 - Senegal
$ws.Range("B80").Value = 2189
$ws.Range("C80").Value = 84
$ws.Range("D80").Value = 842
$ws.Range("E80").Value = 1324
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 23

# Row 130 - Sierra Leona
$ws.Range("B130").Value = 408
$ws.Range("C130").Value = 21
$ws.Range("E130").Value = 285

# Row 163 - Mozambique
$ws.Range("B163").Value = 107
$ws.Range("C163").Value = 3
$ws.Range("D163").Value = 35
$ws.Range("E163").Value = 72
